# Updating filtered feeds from workflow
# Adds two new rows (28 and 29) to the "Filtered Feeds" sheet for the
# GenomeWeb / 360Dx article about Pharma/Dx leaders urging the MRD
# community to focus on clinical utility.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$genomewebUrl = "https://www.genomeweb.com/cancer/pharma-dx-leaders-urge-mrd-community-focus-clinical-utility-use-trials-grows"
$dxUrl        = "https://www.360dx.com/cancer/pharma-dx-leaders-urge-mrd-community-focus-clinical-utility-use-trials-grows"
$keyword      = "MRD"
$title        = "Pharma, Dx Leaders Urge MRD Community to Focus on Clinical Utility as Use in Trials Grows"

# Row 28: GenomeWeb link
$ws.Range("A28").Value = $genomewebUrl
$ws.Range("B28").Value = $keyword
$ws.Range("C28").Value = $title

# Row 29: 360Dx link (same keyword/title, different source URL)
$ws.Range("A29").Value = $dxUrl
$ws.Range("B29").Value = $keyword
$ws.Range("C29").Value = $title

# Register real hyperlinks (and relationships) for the new link cells,
# mirroring the existing rows.
$ws.Hyperlinks.Add($ws.Range("A28"), $genomewebUrl)
$ws.Hyperlinks.Add($ws.Range("A29"), $dxUrl)

# Re-apply the workbook's "Hyperlink" cell style so the new link cells
# match the look (and underlying style index) of the existing ones;
# Hyperlinks.Add() applies its own style, so restore the shared one.
$ws.Range("A28:A29").Style = "Hyperlink"

Write-Host "Added rows 28-29 to Filtered Feeds"
